$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 20240415
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 6
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 8
$ws.Range("G5").Value = 9

$ws.Range("A6").Select()
